$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New text for B4 (FilesTab / dbExcel query): dropped the trailing
# "Study Code" line and the comma after "Diagnosis".
# ---------------------------------------------------------------------
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Rottweiler']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis 
'@

# ---------------------------------------------------------------------
# New text for B2 (CasesTab / dbExcel query): added the Cohort column.
# ---------------------------------------------------------------------
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Rottweiler']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# ---------------------------------------------------------------------
# New text shared by C2/C3/C4 (StatQuery column): replaces the old huge
# filter/aliquot query with a compact program/study/case/sample/file
# roll-up query.
# ---------------------------------------------------------------------
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed  IN ['Rottweiler'] RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Here-strings keep a trailing newline before the closing '@ marker -
# strip that so the text matches the source (no trailing blank line).
$filesQuery = $filesQuery.TrimEnd("`r", "`n")
$casesQuery = $casesQuery.TrimEnd("`r", "`n")
$statQuery  = $statQuery.TrimEnd("`r", "`n")

# Write order matters for shared-string table placement: FilesTab query
# first (new index 11), then the shared StatQuery text (new index 12),
# then the CasesTab+Cohort query last (new index 13) - matching the
# order these strings were (re)introduced upstream.
$ws.Range("B4").Value = $filesQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery
$ws.Range("B2").Value = $casesQuery

# ---------------------------------------------------------------------
# Row heights shrank now that the StatQuery text is much shorter (it no
# longer hits the 409.6pt max-autofit ceiling).
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# ---------------------------------------------------------------------
# View state: selection moved back to B2, no more frozen/scrolled
# topLeftCell at A4.
# ---------------------------------------------------------------------
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
